$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.729.86'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '2.085.76'
$ws.Range("E3").Value = '  +1.28%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''232.43'
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").Value = '''0.623'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '''57.59'
$ws.Range("E8").Value = '  +1.78%  '
$ws.Range("D9").Value = '''0.387'
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("E10").Value = '  +2.34%  '
$ws.Range("E11").Value = '  +3.12%  '
$ws.Range("D12").Value = '2.381.53'
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("D13").Value = '''14.42'
$ws.Range("E13").Value = '  -1.38%  '
$ws.Range("D14").Value = '''21.07'
$ws.Range("D15").Value = '''0.765'
$ws.Range("E15").Value = '  -1.13%  '
$ws.Range("D16").Value = '''5.22'
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("D17").Value = '2.085.52'
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("D18").Value = '37.603.85'
$ws.Range("E18").Value = '  +1.08%  '
$ws.Range("D19").Value = '''6.11'
$ws.Range("E19").Value = '  -3.23%  '
$ws.Range("D20").Value = '''70.61'
$ws.Range("E20").Value = '  +2.16%  '
$ws.Range("E21").Value = '  +1.65%  '
$ws.Range("D22").Value = '''227.87'
$ws.Range("E22").Value = '  +0.98%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").Value = '''2.39'
$ws.Range("E24").Value = '  -1.61%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '''168.08'
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("E27").Value = '  +10.10%  '
$ws.Range("D28").Value = '''8.90'
$ws.Range("E28").Value = '  +1.91%  '
$ws.Range("E29").Value = '  -0.78%  '
$ws.Range("D30").Value = '''19.44'
$ws.Range("E30").Value = '  +2.34%  '
$ws.Range("E31").Value = '  +1.45%  '
$ws.Range("D32").Value = '''4.60'
$ws.Range("E32").Value = '  +4.01%  '
$ws.Range("D33").Value = '''0.0624'
$ws.Range("E33").Value = '  +1.56%  '
$ws.Range("D34").Value = '''4.56'
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("D35").Value = '''2.49'
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("E36").Value = '  +4.05%  '
$ws.Range("D37").Value = '''3.39'
$ws.Range("E37").Value = '  +5.23%  '
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("E39").Value = '  -5.07%  '
$ws.Range("D40").Value = '''0.0997'
$ws.Range("E40").Value = '  +7.38%  '
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("D42").Value = '''97.58'
$ws.Range("E42").Value = '  +1.65%  '
$ws.Range("D43").Value = '''0.0213'
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("D44").Value = '1.452.40'
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("E45").Value = '  -0.55%  '
$ws.Range("D46").Value = '''1.05'
$ws.Range("E46").Value = '  +3.51%  '
$ws.Range("D47").Value = '''4.07'
$ws.Range("E47").Value = '  -4.07%  '
$ws.Range("D48").Value = '''15.64'
$ws.Range("D49").Value = '''7.33'
$ws.Range("E49").Value = '  +2.80%  '
$ws.Range("D50").Value = '''3.00'
$ws.Range("E50").Value = '  +1.81%  '
$ws.Range("D51").Value = '2.276.87'
$ws.Range("E51").Value = '  +1.25%  '
